$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")
$ws.Range("C4").Value = 750
